# Add fields to current_inventory table:
# Replace rows 3-6 (Heidtman Cold Rolled Steel Sheet data) with the last batch
# of Steel Warehouse - 93158 Hot Rolled Steel Sheet records (TCY8024/8026/8025/8027),
# and remove row 7 (the data set shrank from 5 rows to 4 rows for this block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 3 (TCY8024) ----
$ws.Range("B3").Value = "Steel Warehouse - 93158"
$ws.Range("D3").Value = "Hot Rolled Steel Sheet"
$ws.Range("K3").Value = "TCY8024"
$ws.Range("M3").Value = "2759J5 02"
$ws.Range("O3").Value = "1200025024 / 2203158 / L02524/1"
$ws.Range("S3").Value = "0.0580"
$ws.Range("V3").Value = "48.500"
$ws.Range("AF3").Value = "45,200"

# ---- Row 4 (TCY8026) ----
$ws.Range("B4").Value = "Steel Warehouse - 93158"
$ws.Range("D4").Value = "Hot Rolled Steel Sheet"
$ws.Range("K4").Value = "TCY8026"
$ws.Range("M4").Value = "2759J5 03"
$ws.Range("O4").Value = "1200025024 / 2203158 / L02524/1"
$ws.Range("S4").Value = "0.0580"
$ws.Range("V4").Value = "48.500"
$ws.Range("AF4").Value = "46,220"

# ---- Row 5 (TCY8025) ----
$ws.Range("B5").Value = "Steel Warehouse - 93158"
$ws.Range("D5").Value = "Hot Rolled Steel Sheet"
$ws.Range("K5").Value = "TCY8025"
$ws.Range("M5").Value = "2759J5 52"
$ws.Range("O5").Value = "1200025024 / 2203158 / L02524/1"
$ws.Range("S5").Value = "0.0580"
$ws.Range("V5").Value = "48.500"
$ws.Range("AF5").Value = "44,860"

# ---- Row 6 (TCY8027) ----
$ws.Range("B6").Value = "Steel Warehouse - 93158"
$ws.Range("D6").Value = "Hot Rolled Steel Sheet"
$ws.Range("K6").Value = "TCY8027"
$ws.Range("M6").Value = "2759J5 53"
$ws.Range("O6").Value = "1200025024 / 2203158 / L02524/1"
$ws.Range("S6").Value = "0.0580"
$ws.Range("V6").Value = "48.500"
$ws.Range("AF6").Value = "43,820"

# The numeric-looking values above (S, V, AF columns) must stay stored as literal
# text (matching the source shared-string formatting, e.g. "0.0580" and "45,200")
# instead of being auto-converted to numbers. Force text format before assigning,
# then clear the style stamp that NumberFormat="@" leaves behind so the cells keep
# their original (unstyled) appearance.
$numericRanges = @("S3:S6", "V3:V6", "AF3:AF6")
foreach ($rng in $numericRanges) {
    $ws.Range($rng).NumberFormat = "@"
}

$ws.Range("S3").Value = "0.0580"
$ws.Range("S4").Value = "0.0580"
$ws.Range("S5").Value = "0.0580"
$ws.Range("S6").Value = "0.0580"

$ws.Range("V3").Value = "48.500"
$ws.Range("V4").Value = "48.500"
$ws.Range("V5").Value = "48.500"
$ws.Range("V6").Value = "48.500"

$ws.Range("AF3").Value = "45,200"
$ws.Range("AF4").Value = "46,220"
$ws.Range("AF5").Value = "44,860"
$ws.Range("AF6").Value = "43,820"

foreach ($rng in $numericRanges) {
    $ws.Range($rng).Style = "Normal"
}

# Row 7 no longer exists in the updated data set - remove it entirely, which
# also shrinks the used range from A1:AS7 down to A1:AS6.
$ws.Rows.Item(7).Delete()

Write-Host "Applied current_inventory field updates"
